$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new row 67's date cell inherits the same date format as the other date cells in column D
$ws.Range("D67").NumberFormat = $ws.Range("D2").NumberFormat

$data = @(
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44162, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 85, 2200, 2300, 2247, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2247, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44406, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 100, 3200, 3200, 3200, '$/kilo (en caja de 15 kilos)', 'Provincia de Limarí', 3200, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44455, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 15, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44455, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 20, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44442, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 35, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44442, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 40, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44354, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Tercera', 95, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44459, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 80, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44459, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 50, 2800, 2800, 2800, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2800, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44420, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 35, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44420, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 40, 3200, 3200, 3200, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3200, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44160, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 120, 2200, 2300, 2246, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2246, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44431, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 80, 21000, 21000, 21000, '$/bandeja 7 kilos', 'Provincia del Elquí', 3000, 7),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44460, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 20, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44460, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 60, 2800, 3000, 2900, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2900, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44405, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 50, 3200, 3200, 3200, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3200, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44446, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 30, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44446, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 80, 2800, 2800, 2800, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2800, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44417, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 50, 3200, 3200, 3200, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3200, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44419, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 70, 3200, 3200, 3200, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3200, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44424, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 25, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44449, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 50, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44412, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 65, 3200, 3200, 3200, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3200, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44428, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 55, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44447, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 40, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44421, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 50, 3200, 3200, 3200, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3200, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44454, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 55, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44454, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 35, 3200, 3200, 3200, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3200, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44467, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Sin especificar', 'Primera', 30, 3200, 3200, 3200, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3200, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44427, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 65, 24000, 24000, 24000, '$/bandeja 7 kilos', 'Provincia del Elquí', 3429, 7),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44473, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 50, 22000, 22000, 22000, '$/bandeja 7 kilos', 'Provincia del Elquí', 3143, 7),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44411, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 10, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44413, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 35, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44426, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 30, 4500, 4500, 4500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 4500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44426, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 45, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44469, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 90, 3000, 3200, 3111, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3111, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44434, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 80, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44161, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 65, 2300, 2300, 2300, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2300, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44161, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 55, 2000, 2000, 2000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44407, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 40, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44165, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 50, 2300, 2300, 2300, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2300, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44475, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 20, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44474, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 20, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44474, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 60, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44448, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 40, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia de Limarí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44448, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 200, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia de Limarí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44452, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 45, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44452, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 65, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44453, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 65, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44453, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 55, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44435, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 130, 3500, 4500, 3885, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3885, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44435, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 80, 21000, 21000, 21000, '$/bandeja 7 kilos', 'Provincia del Elquí', 3000, 7),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44377, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 40, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44433, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 20, 4500, 4500, 4500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 4500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44159, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 120, 2300, 2500, 2408, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2408, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44414, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 55, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44425, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 35, 4500, 4500, 4500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 4500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44425, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 20, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44425, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 25, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44187, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 15, 3200, 3200, 3200, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3200, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44438, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 35, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44438, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 20, 3000, 3000, 3000, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3000, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44461, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 40, 3500, 3500, 3500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 3500, 1),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44461, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 40, 30000, 30000, 30000, '$/bandeja 8 kilos', 'Provincia del Elquí', 3750, 8),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44461, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Segunda', 30, 28000, 28000, 28000, '$/bandeja 8 kilos', 'Provincia del Elquí', 3500, 8),
    @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44432, 9, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 30, 4500, 4500, 4500, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 4500, 1)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

$ws.Range("A1").Select() | Out-Null
